$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 7) with the same shape as the existing "liens" rows:
# col A = "TEST" (category marker, same as every other row),
# col B = label, col C = description, col D = a hyperlink cell.
$ws.Range("A7").Value = "TEST"
$ws.Range("B7").Value = "dff"
$ws.Range("C7").Value = "dfsfdsfsd"

# D7 becomes a real hyperlink to https://youtube.com, displaying the URL text
# (Excel auto-fills the cell text with the address when no TextToDisplay is given).
$ws.Hyperlinks.Add($ws.Range("D7"), "https://youtube.com")

# Match the formatting used by the rest of the table: copy the (unstyled/default)
# look of row 6's A:C cells onto the new row 7's A:C cells.
$ws.Range("A6:C6").Copy()
$ws.Range("A7:C7").PasteSpecial(-4122)  # xlPasteFormats
